# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows above the current row 2 (existing data shifts down by 6 rows)
$ws.Rows("2:7").Insert()
$ws.Rows("2:7").ClearFormats()

# New data for the inserted rows (new rows 2-7)
$topData = @(
    @(-0.0641408488154411, 0.0114537235349416, -0.641255795955658),
    @(-0.0258090570569038, -0.0752891451120376, -0.038026362657547),
    @(-0.2993239760398865, 0.2926044464111328, 0.0584903471171855),
    @(-0.08659014850854869, -0.0867428630590438, 0.1705841124057769),
    @(0.1209513172507286, -0.2257147133350372, -0.1018617823719978),
    @(-0.070249505341053, -0.1577559560537338, 0.0288633834570646)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $topData[$i][0]
    $ws.Cells.Item($row, 2).Value = $topData[$i][1]
    $ws.Cells.Item($row, 3).Value = $topData[$i][2]
}

# New data appended at the bottom (new rows 28-31)
$bottomData = @(
    @(-0.4198171496391296, 0.1876883506774902, 0.732122004032135),
    @(0.0713185146450996, -0.2884811162948608, 0.1009454801678657),
    @(0.0609338097274303, -0.087353728711605, 0.0229074470698833),
    @(0.015118914656341, -0.0455094613134861, 0.119576871395111)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $row = 28 + $i
    $ws.Cells.Item($row, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($row, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($row, 3).Value = $bottomData[$i][2]
}
